# Apply the "DTR summary" edits:
#  - remove the excess "1" values that had leaked into the I column on a
#    few daily rows (these were stray numbers, not meant to be there)
#  - add a "Legends:" section under the summary table explaining the
#    color-coded remarks used elsewhere in the report

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove excess/erroneous values in column I (rows 9, 10, 14) -------
$ws.Range("I9").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("I14").ClearContents()

# --- 2. Legends section ----------------------------------------------------

# "Legends:" heading, styled like the big report title (Arial 15 bold underline)
$ws.Range("E24:P24").Merge()
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null          # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E24").Value = "Legends:"

# Blue swatch + explanation: requests/remarks for the day
$ws.Range("E25:E26").Merge()
$ws.Range("F25:P26").Merge()
$ws.Range("E25:E26").Interior.Color = 13411113            # FF29A3CC
$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."

# Orange swatch + explanation: half-day
$ws.Range("E27:E28").Merge()
$ws.Range("F27:P28").Merge()
$ws.Range("E27:E28").Interior.Color = 6737151              # FFFFCC66
$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."

# Red swatch + explanation: absent
$ws.Range("E29:E30").Merge()
$ws.Range("F29:P30").Merge()
$ws.Range("E29:E30").Interior.Color = 6184671               # FFDF5E5E
$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
